$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking text like "1.016" or "0.00001112".
# Assigning such a string directly would make Excel reinterpret it as a real
# number (and e.g. drop trailing zeros / switch to scientific notation), so we
# temporarily force the range to Text format while writing the values, then
# restore the default "Normal" style so no stray number-format style sticks
# around on the cells (matches the original workbook, which has no explicit
# style on these cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.422.83"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.895.88"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("D4").Value = "1.016"
$ws.Range("E4").Value = "  +1.19%  "
$ws.Range("D5").Value = "316.98"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").Value = "1.014"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("D7").Value = "0.5176"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("D8").Value = "0.3942"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "0.08443"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").Value = "1.130"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "42.01"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "6.303"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "20.77"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "1.880.83"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "7.337"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "1.016"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "0.00001112"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "91.86"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "0.06748"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").Value = "17.94"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "1.014"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").Value = "6.083"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "28.506.01"
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("D24").Value = "11.23"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "2.272"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "161.04"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.492"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "20.93"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").Value = "127.27"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").Value = "0.1062"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "1.045"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").Value = "5.867"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "3.640"
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").Value = "9.718"
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").Value = "0.02477"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("D36").Value = "0.06629"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").Value = "0.2224"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").Value = "1.209"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "0.6544"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").Value = "1.246"
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").Value = "5.013"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").Value = "11.37"
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").Value = "0.6160"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").Value = "13.12"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "3.707"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.289"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "2.033"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").Value = "1.244"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").Value = "121.60"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("D50").Value = "0.06952"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("D51").Value = "78.34"
$ws.Range("E51").Value = "  +0.14%  "

# Drop the temporary Text number format back to the workbook default so the
# cells end up with no explicit style, same as before the edit.
$ws.Range("D2:D51").Style = "Normal"
